$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped from 45171 to
# 45172 (i.e. 2023-09-02 -> 2023-09-03) for every data row (rows 2-246).
$ws.Range("C2:C246").Value = 45172
